# Add hydrogen emissions data (year 2050) to the "year_Vecteurs" sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("year_Vecteurs")

# Update existing row 7 (hydrogen / 2020): indirect_emissions (D7) 0 -> 289
$ws.Range("D7").Value = 289

# Append new rows 8-13 for year 2050, mirroring rows 2-7 (elec, gaz, fioul, bois, charbon, hydrogen)
$ws.Range("A8").Value = "elec"
$ws.Range("B8").Value = 2050
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 15

$ws.Range("A9").Value = "gaz"
$ws.Range("B9").Value = 2050
$ws.Range("C9").Value = 44
$ws.Range("D9").Value = 0

$ws.Range("A10").Value = "fioul"
$ws.Range("B10").Value = 2050
$ws.Range("C10").Value = 272
$ws.Range("D10").Value = 57

$ws.Range("A11").Value = "bois"
$ws.Range("B11").Value = 2050
$ws.Range("C11").Value = 27
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "charbon"
$ws.Range("B12").Value = 2050
$ws.Range("C12").Value = 346.5
$ws.Range("D12").Value = 28.5

$ws.Range("A13").Value = "hydrogen"
$ws.Range("B13").Value = 2050
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 15

# Match the highlighted (yellow) style used on the hydrogen row (C7:D7) for the new hydrogen row
$ws.Range("C7:D7").Copy()
$ws.Range("C13:D13").PasteSpecial(-4122) # xlPasteFormats

# Make this sheet the active tab, with the selection on G6 as in the target file
$ws.Range("G6").Select()
$ws.Activate()
